# Update scripts with new TPM values (Edn3-Ednrb LR pair table)
# - Adds "ECs" as a new sending cluster (rows 2-5), shifting the former
#   FAPs/MuSCs sending-cluster rows down, and appends the corresponding
#   MuSCs sending-cluster rows (10-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Edn3"
$ws.Cells.Item(2, 3).Value = "Ednrb"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1264883333333333
$ws.Cells.Item(2, 8).Value = 0.379465
$ws.Cells.Item(2, 9).Value = 0.02088586470611676
$ws.Cells.Item(2, 10).Value = 0.02088586470611676
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 36.27867833333334
$ws.Cells.Item(2, 14).Value = 108.836035
$ws.Cells.Item(2, 15).Value = 0.6141201546381995
$ws.Cells.Item(2, 16).Value = 0.6141201546381995
$ws.Cells.Item(2, 17).Value = 4.588829557919445
$ws.Cells.Item(2, 18).Value = 41.299466021275
$ws.Cells.Item(2, 19).Value = 0.01282643046307294
$ws.Cells.Item(2, 20).Value = 0.01282643046307294

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Edn3"
$ws.Cells.Item(3, 3).Value = "Ednrb"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1264883333333333
$ws.Cells.Item(3, 8).Value = 0.379465
$ws.Cells.Item(3, 9).Value = 0.02088586470611676
$ws.Cells.Item(3, 10).Value = 0.02088586470611676
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.152389
$ws.Cells.Item(3, 14).Value = 0.457167
$ws.Cells.Item(3, 15).Value = 0.002579618678092064
$ws.Cells.Item(3, 16).Value = 0.002579618678092065
$ws.Cells.Item(3, 17).Value = 0.01927543062833333
$ws.Cells.Item(3, 18).Value = 0.173478875655
$ws.Cells.Item(3, 19).Value = 0.00005387756670400262
$ws.Cells.Item(3, 20).Value = 0.00005387756670400263

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Edn3"
$ws.Cells.Item(4, 3).Value = "Ednrb"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1264883333333333
$ws.Cells.Item(4, 8).Value = 0.379465
$ws.Cells.Item(4, 9).Value = 0.02088586470611676
$ws.Cells.Item(4, 10).Value = 0.02088586470611676
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 20.63960233333333
$ws.Cells.Item(4, 14).Value = 61.918807
$ws.Cells.Item(4, 15).Value = 0.3493841660976791
$ws.Cells.Item(4, 16).Value = 0.3493841660976791
$ws.Cells.Item(4, 17).Value = 2.610668899806111
$ws.Cells.Item(4, 18).Value = 23.496020098255
$ws.Cells.Item(4, 19).Value = 0.007297190423575553
$ws.Cells.Item(4, 20).Value = 0.007297190423575553

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Edn3"
$ws.Cells.Item(5, 3).Value = "Ednrb"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1264883333333333
$ws.Cells.Item(5, 8).Value = 0.379465
$ws.Cells.Item(5, 9).Value = 0.02088586470611676
$ws.Cells.Item(5, 10).Value = 0.02088586470611676
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.003565333333333
$ws.Cells.Item(5, 14).Value = 6.010696
$ws.Cells.Item(5, 15).Value = 0.03391606058602931
$ws.Cells.Item(5, 16).Value = 0.03391606058602931
$ws.Cells.Item(5, 17).Value = 0.2534276397377778
$ws.Cells.Item(5, 18).Value = 2.28084875764
$ws.Cells.Item(5, 19).Value = 0.0007083662527642672
$ws.Cells.Item(5, 20).Value = 0.0007083662527642673

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Edn3"
$ws.Cells.Item(6, 3).Value = "Ednrb"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.012975
$ws.Cells.Item(6, 8).Value = 0.038925
$ws.Cells.Item(6, 9).Value = 0.002142443397113291
$ws.Cells.Item(6, 10).Value = 0.002142443397113291
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 36.27867833333334
$ws.Cells.Item(6, 14).Value = 108.836035
$ws.Cells.Item(6, 15).Value = 0.6141201546381995
$ws.Cells.Item(6, 16).Value = 0.6141201546381995
$ws.Cells.Item(6, 17).Value = 0.4707158513750001
$ws.Cells.Item(6, 18).Value = 4.236442662375
$ws.Cells.Item(6, 19).Value = 0.001315717670338803
$ws.Cells.Item(6, 20).Value = 0.001315717670338804

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Edn3"
$ws.Cells.Item(7, 3).Value = "Ednrb"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.012975
$ws.Cells.Item(7, 8).Value = 0.038925
$ws.Cells.Item(7, 9).Value = 0.002142443397113291
$ws.Cells.Item(7, 10).Value = 0.002142443397113291
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.152389
$ws.Cells.Item(7, 14).Value = 0.457167
$ws.Cells.Item(7, 15).Value = 0.002579618678092064
$ws.Cells.Item(7, 16).Value = 0.002579618678092065
$ws.Cells.Item(7, 17).Value = 0.001977247275
$ws.Cells.Item(7, 18).Value = 0.017795225475
$ws.Cells.Item(7, 19).Value = 0.000005526687003948459
$ws.Cells.Item(7, 20).Value = 0.00000552668700394846

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Edn3"
$ws.Cells.Item(8, 3).Value = "Ednrb"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.012975
$ws.Cells.Item(8, 8).Value = 0.038925
$ws.Cells.Item(8, 9).Value = 0.002142443397113291
$ws.Cells.Item(8, 10).Value = 0.002142443397113291
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 20.63960233333333
$ws.Cells.Item(8, 14).Value = 61.918807
$ws.Cells.Item(8, 15).Value = 0.3493841660976791
$ws.Cells.Item(8, 16).Value = 0.3493841660976791
$ws.Cells.Item(8, 17).Value = 0.267798840275
$ws.Cells.Item(8, 18).Value = 2.410189562475
$ws.Cells.Item(8, 19).Value = 0.0007485357997119059
$ws.Cells.Item(8, 20).Value = 0.000748535799711906

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Edn3"
$ws.Cells.Item(9, 3).Value = "Ednrb"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.012975
$ws.Cells.Item(9, 8).Value = 0.038925
$ws.Cells.Item(9, 9).Value = 0.002142443397113291
$ws.Cells.Item(9, 10).Value = 0.002142443397113291
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.003565333333333
$ws.Cells.Item(9, 14).Value = 6.010696
$ws.Cells.Item(9, 15).Value = 0.03391606058602931
$ws.Cells.Item(9, 16).Value = 0.03391606058602931
$ws.Cells.Item(9, 17).Value = 0.0259962602
$ws.Cells.Item(9, 18).Value = 0.2339663418
$ws.Cells.Item(9, 19).Value = 0.00007266324005863281
$ws.Cells.Item(9, 20).Value = 0.00007266324005863284

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Edn3"
$ws.Cells.Item(10, 3).Value = "Ednrb"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.916706
$ws.Cells.Item(10, 8).Value = 17.750118
$ws.Cells.Item(10, 9).Value = 0.9769716918967699
$ws.Cells.Item(10, 10).Value = 0.97697169189677
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 36.27867833333334
$ws.Cells.Item(10, 14).Value = 108.836035
$ws.Cells.Item(10, 15).Value = 0.6141201546381995
$ws.Cells.Item(10, 16).Value = 0.6141201546381995
$ws.Cells.Item(10, 17).Value = 214.6502737669034
$ws.Cells.Item(10, 18).Value = 1931.85246390213
$ws.Cells.Item(10, 19).Value = 0.5999780065047877
$ws.Cells.Item(10, 20).Value = 0.5999780065047877

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Edn3"
$ws.Cells.Item(11, 3).Value = "Ednrb"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 5.916706
$ws.Cells.Item(11, 8).Value = 17.750118
$ws.Cells.Item(11, 9).Value = 0.9769716918967699
$ws.Cells.Item(11, 10).Value = 0.97697169189677
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.152389
$ws.Cells.Item(11, 14).Value = 0.457167
$ws.Cells.Item(11, 15).Value = 0.002579618678092064
$ws.Cells.Item(11, 16).Value = 0.002579618678092065
$ws.Cells.Item(11, 17).Value = 0.9016409106340001
$ws.Cells.Item(11, 18).Value = 8.114768195706
$ws.Cells.Item(11, 19).Value = 0.002520214424384113
$ws.Cells.Item(11, 20).Value = 0.002520214424384114

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Edn3"
$ws.Cells.Item(12, 3).Value = "Ednrb"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 5.916706
$ws.Cells.Item(12, 8).Value = 17.750118
$ws.Cells.Item(12, 9).Value = 0.9769716918967699
$ws.Cells.Item(12, 10).Value = 0.97697169189677
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 20.63960233333333
$ws.Cells.Item(12, 14).Value = 61.918807
$ws.Cells.Item(12, 15).Value = 0.3493841660976791
$ws.Cells.Item(12, 16).Value = 0.3493841660976791
$ws.Cells.Item(12, 17).Value = 122.1184589632473
$ws.Cells.Item(12, 18).Value = 1099.066130669226
$ws.Cells.Item(12, 19).Value = 0.3413384398743917
$ws.Cells.Item(12, 20).Value = 0.3413384398743917

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Edn3"
$ws.Cells.Item(13, 3).Value = "Ednrb"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 5.916706
$ws.Cells.Item(13, 8).Value = 17.750118
$ws.Cells.Item(13, 9).Value = 0.9769716918967699
$ws.Cells.Item(13, 10).Value = 0.97697169189677
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.003565333333333
$ws.Cells.Item(13, 14).Value = 6.010696
$ws.Cells.Item(13, 15).Value = 0.03391606058602931
$ws.Cells.Item(13, 16).Value = 0.03391606058602931
$ws.Cells.Item(13, 17).Value = 11.85450702912533
$ws.Cells.Item(13, 18).Value = 106.690563262128
$ws.Cells.Item(13, 19).Value = 0.0331350310932064
$ws.Cells.Item(13, 20).Value = 0.03313503109320642

